# Change for button template
$wb = $excel.ActiveWorkbook

$wsCards = $wb.Worksheets.Item("RICH_CARDS_CONFIG")

# Update the two QUICK_REPLIES rows in RICH_CARDS_CONFIG to BUTTONS
$wsCards.Range("B2").Value = "BUTTONS"
$wsCards.Range("B3").Value = "BUTTONS"

# Move the selection/active cell on RICH_CARDS_CONFIG to B4
$wsCards.Range("B4").Select()

# Make RICH_CARDS_CONFIG the active (selected) tab instead of PROMPTS_CONFIG
$wsCards.Activate()
